$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card3")

# N1 header: remove trailing space -> "Event"
$ws.Range("N1").Value = "Event"

# N2-N7: fill with "nan"
$ws.Range("N2").Value = "nan"
$ws.Range("N3").Value = "nan"
$ws.Range("N4").Value = "nan"
$ws.Range("N5").Value = "nan"
$ws.Range("N6").Value = "nan"
$ws.Range("N7").Value = "nan"

# N8: Arabic service note
$ws.Range("N8").Value = "تم سن الفلاتس لاول مره بعد التغير"

# N9-N13: fill with "nan"
$ws.Range("N9").Value = "nan"
$ws.Range("N10").Value = "nan"
$ws.Range("N11").Value = "nan"
$ws.Range("N12").Value = "nan"
$ws.Range("N13").Value = "nan"
